$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Updated monomeric yield values (T and Y columns) for rows 4-12, and
# new (blank, formatted) helper cells in columns U and Z alongside them.
$rows = @(
    @{ Row = 4;  T = 4.9775632205666147; Y = 0.42108293339177427 },
    @{ Row = 5;  T = 5.1077837208470473; Y = 0.55130343367220691 },
    @{ Row = 6;  T = 6.5387604656145548; Y = 0.99968339351528279 },
    @{ Row = 7;  T = 6.8533893643976977; Y = 1.3825615637196409 },
    @{ Row = 8;  T = 7.0950106164539717; Y = 1.8496484107134943 },
    @{ Row = 9;  T = 5.6416721059429582; Y = 5.601799697649497 },
    @{ Row = 10; T = 5.5145464936056126; Y = 4.2386294282148631 },
    @{ Row = 11; T = 4.8163851955016499; Y = 1.5402023140556043 },
    @{ Row = 12; T = 3.6515751376528556; Y = 0.27571123547316245 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("T$rowNum").Value = $r.T
    $ws.Range("Y$rowNum").Value = $r.Y
    # Touch the new neighbouring cells so they exist in the sheet (blank,
    # but carrying a style) just like in the authored workbook.
    $ws.Range("U$rowNum").Value = ""
    $ws.Range("Z$rowNum").Value = ""
}

# Remove the now-superseded "120 min" isothermal-time data row; this
# shifts every row below it (the helper grid starting at row 17, and
# the trailing blank formatted rows) up by one.
$ws.Rows.Item(13).Delete()

# Reflect the final cell selection left by the author.
$ws.Range("J23").Select()
